$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ58436446"
$wb.Worksheets.Item(2).Name = "summ58555268"
$wb.Worksheets.Item(3).Name = "summ58684597"
$wb.Worksheets.Item(4).Name = "summ58812334"
$wb.Worksheets.Item(5).Name = "summ58933402"
$wb.Worksheets.Item(6).Name = "summ59046852"
$wb.Worksheets.Item(7).Name = "summ59157323"
$wb.Worksheets.Item(8).Name = "summ59354362"
$wb.Worksheets.Item(9).Name = "summ59451383"
$wb.Worksheets.Item(10).Name = "summ59549702"
$wb.Worksheets.Item(11).Name = "summ59646695"
$wb.Worksheets.Item(12).Name = "summ59743933"
$wb.Worksheets.Item(13).Name = "summ59841867"
$wb.Worksheets.Item(14).Name = "summ59940882"
$wb.Worksheets.Item(15).Name = "summ00062406"
$wb.Worksheets.Item(16).Name = "summ00176307"
$wb.Worksheets.Item(17).Name = "summ00284299"
$wb.Worksheets.Item(18).Name = "summ00405540"
$wb.Worksheets.Item(19).Name = "summ00531232"
$wb.Worksheets.Item(20).Name = "summ00658031"
$wb.Worksheets.Item(21).Name = "summ00782242"
$wb.Worksheets.Item(22).Name = "summ00911550"
$wb.Worksheets.Item(23).Name = "summ01038364"
$wb.Worksheets.Item(24).Name = "summ01164645"
$wb.Worksheets.Item(25).Name = "summ01295563"
$wb.Worksheets.Item(26).Name = "summ01427734"
$wb.Worksheets.Item(27).Name = "summ01560108"
$wb.Worksheets.Item(28).Name = "summ01724463"
$wb.Worksheets.Item(29).Name = "summ01888768"
$wb.Worksheets.Item(30).Name = "summ02031171"
$wb.Worksheets.Item(31).Name = "summ02161232"
$wb.Worksheets.Item(32).Name = "summ02299729"
$wb.Worksheets.Item(33).Name = "summ02432067"
$wb.Worksheets.Item(34).Name = "summ02555534"
$wb.Worksheets.Item(35).Name = "summ02680907"
$wb.Worksheets.Item(36).Name = "summ02834245"
$wb.Worksheets.Item(37).Name = "summ02966649"
$wb.Worksheets.Item(38).Name = "summ03093168"
$wb.Worksheets.Item(39).Name = "summ03214071"
$wb.Worksheets.Item(40).Name = "summ03334814"
$wb.Worksheets.Item(41).Name = "summ03455065"
$wb.Worksheets.Item(42).Name = "summ03586880"
$wb.Worksheets.Item(43).Name = "summ03730216"
$wb.Worksheets.Item(44).Name = "summ03863386"
$wb.Worksheets.Item(45).Name = "summ03983117"
$wb.Worksheets.Item(46).Name = "summ04108102"
$wb.Worksheets.Item(47).Name = "summ04251940"
$wb.Worksheets.Item(48).Name = "summ04381148"
$wb.Worksheets.Item(49).Name = "summ04518822"
$wb.Worksheets.Item(50).Name = "summ04662350"
